# "Redigerte en feil i ukedagene" (Fixed an error in the weekdays)
#
# The hours that were entered under "Dag 8" (column L) actually belonged to
# "Dag 9" (column M), so the existing column-L numbers are shifted one day
# to the right into column M, and column L is refilled by repeating the
# previous day's ("Dag 7" / column K) figure - i.e. no extra work recorded
# that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-18, 20 and 21 currently have an M cell whose formatting (style)
# differs from column L's; use Copy so the destination cell picks up both
# the value AND the exact formatting of the source cell.
$copyRows = @(3, 10, 11, 12, 13, 14, 15, 16, 17, 18, 20, 21)
foreach ($r in $copyRows) {
    $ws.Range("L$r").Copy($ws.Range("M$r"))
}

# Remaining rows already carry the correct destination formatting, so just
# move the value itself.
$valueOnlyRows = @(2, 4, 5, 6, 7, 8, 9, 19)
foreach ($r in $valueOnlyRows) {
    $ws.Range("M$r").Value = $ws.Range("L$r").Value2
}

# Column L is repopulated with the previous day's (column K) figure for
# every data row (2-21).
foreach ($r in 2..21) {
    $ws.Range("L$r").Value = $ws.Range("K$r").Value2
}

# Restore the last-clicked cell recorded in the saved file.
$ws.Range("O7").Select() | Out-Null
